$d = $word.ActiveDocument

# 1. Remove the narrative about the motion to amend the charge (and the
#    Court's ruling granting it) from the arraignment paragraph, leaving
#    the trailing empty run that followed "for arraignment on January 07, 2022. "
$removedText = "Counsel for the State of Ohio made a motion to amend the charge of DUS UCM to Assured Clear Distrance Ahead. The Court found the amendment did not alter the name or identity of the offense and the motion is Granted. "
$d.Content.Find.Execute($removedText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2) | Out-Null

# 2. The charge description in the table reverts from "DUS UCM - AMENDED"
#    back to "DUS UCM" (the amendment was not carried through after all).
$tbl = $d.Tables.Item(1)
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
        $cell = $tbl.Cell($r, $c)
        if ($cell.Range.Text -like "*DUS UCM - AMENDED*") {
            $rng = $cell.Range
            $rng.End = $rng.End - 1
            $rng.Text = "DUS UCM"
        }
    }
}
